$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# B7 is the value cell next to "Experimental" (A7) - set it to the literal text
# "true" (not the Excel boolean TRUE). Writing the bare word directly would be
# auto-converted to a boolean by Excel, so build it as a formula result and
# paste back as a value, which preserves it as plain text.
$scratch = $ws.Cells.Item(20, 2)
$scratch.Formula = '="true"'
$scratch.Copy()
$ws.Cells.Item(7, 2).PasteSpecial(-4163)
$scratch.Clear()

# B8 is the value cell next to "Date" (A8) - update the date string
$ws.Range("B8").Value = "2023-02-16T14:43:10-06:00"
